# Partial implementation of rich text for dialogue.
# Dialogue text in cells C4 then C3 is rewritten to use Unity rich-text tags.
# (C4 is set first so the shared-string table ends up ordered the same way
# Excel produced it: the "Ribbit..." string before the "Hey you!..." string.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "Ribbit <i>Ribbit!</i> (Yeah <color=green>frog-face!</color> Wrong part of town!)"
$ws.Range("C3").Value = "<size=48>Hey you!</size> You're walking in the <color=red>wrong</color> part of town."

# Reflect the author's final cell selection in the saved worksheet view.
$ws.Range("C3").Select()
